# Restructure the "healthscore_calculado" sheet:
#  - the "tags" column (old col E) is dropped and the duplicated trailing
#    "Healthscore" column (old col K) is removed, so ocorridos / data / cliente /
#    Ranking_de_Eventos / Delta / Healthscore all shift one column to the left
#    (columns now run A..J instead of A..K)
#  - "data" (now column F) becomes a real date serial value, displayed with a
#    yyyy-mm-dd hh:mm:ss style instead of being a plain text date string
#  - the FBMDS/Litero rows are corrected and reordered, and 3 new backlog rows
#    (Cliente_fake, Mart Minas, Mart Minas) are appended, growing the used
#    range from A1:K10 to A1:J13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: overwrite the text in place so the existing bold / bordered /
#     centered header style already on A1:K1 is kept as-is ---
$headers = @('id Runrunit', 'titulo', 'estado', 'Quadro', 'ocorridos', 'data', 'cliente', 'Ranking_de_Eventos', 'Delta', 'Healthscore')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Old column K ("tags" data lived under E, and K used to hold a second,
# now-redundant copy of "Healthscore") is no longer part of the table.
$ws.Range("K1:K10").Clear()

# --- Data rows ---
# columns: id Runrunit, titulo, estado, Quadro, ocorridos, data, cliente,
#          Ranking_de_Eventos, Delta, Healthscore
# "data" is stored as an Excel date serial number (e.g. 45488 = 2024-07-15).
# "cliente" is intentionally blank for the three anonymized "Cliente_fake" rows.
$data = @(
    @(157821, 'Cocamar, 15/07/2024', 'backlog', 'Acompanhamento de clientes', 'Cliente pediu proposta', 45488, 'Cocamar', 'Cliente pediu proposta', 2.5, 10),
    @(157821, 'Cocamar, 15/07/2024', 'backlog', 'Acompanhamento de clientes', 'Resolveu problema', 45488, 'Cocamar', 'Resolveu problema', 2, 10),
    @(149896, 'Mart Minas, 15/07/2024', 'backlog', 'Acompanhamento de clientes', 'Metas não atingidas', 45488, 'Mart Minas', 'Metas não atingidas', -2, 8),
    @(149896, 'Mart Minas, 15/07/2024', 'backlog', 'Acompanhamento de clientes', 'Feedback positivo', 45488, 'Mart Minas', 'Feedback positivo', 2.5, 10),
    @(156244, 'FBMDS, 29/07/2024', 'backlog', 'Acompanhamento de clientes', 'Planejamento foi reprovado', 45502, 'FBMDS', 'Planejamento foi reprovado', -1, 9),
    @(156243, 'Litero, 05/08/2024', 'backlog', 'Acompanhamento de clientes', 'Cliente pediu proposta', 45509, 'Litero', 'Cliente pediu proposta', 2.5, 10),
    @(156243, 'Litero, 05/08/2024', 'backlog', 'Acompanhamento de clientes', 'Feedback positivo', 45509, 'Litero', 'Feedback positivo', 2.5, 10),
    @(160421, 'Cliente_fake, 19/08/2024', 'backlog', 'Acompanhamento de clientes', 'Feedback negativo', 45523, '', 'Feedback negativo', -2.5, 7.5),
    @(160421, 'Cliente_fake, 19/08/2024', 'backlog', 'Acompanhamento de clientes', 'Cliente não responde (NPS ou outra comunicação) | No-show', 45523, '', 'Cliente não responde (NPS ou outra comunicação) | No-show', -2, 5.5),
    @(160421, 'Cliente_fake, 19/08/2024', 'backlog', 'Acompanhamento de clientes', 'Planejamento foi reprovado', 45523, '', 'Planejamento foi reprovado', -1, 4.5),
    @(160794, 'Mart Minas, 26/08/2024', 'backlog', 'Acompanhamento de clientes', 'Planejamento foi reprovado', 45530, 'Mart Minas', 'Planejamento foi reprovado', -1, 9),
    @(165463, 'Mart Minas, 09/09/2024', 'backlog', 'Acompanhamento de clientes', 'Feedback negativo', 45544, 'Mart Minas', 'Feedback negativo', -2.5, 6.5)
)

$dateFormatted = $false
$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    # Column 6 ("data") holds a date serial number - give it a date/time display
    # format (registers the yyyy-mm-dd / YYYY-MM-DD number formats, matching
    # the two numFmts the workbook ends up with).
    if (-not $dateFormatted) {
        $ws.Cells.Item($r, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
        $dateFormatted = $true
    }
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r++
}
